# This document currently has no word/styles.xml part at all -- every
# paragraph just falls back to the implicit built-in "Normal" style.
# Materialize that style explicitly by touching the Styles collection,
# which mints a real styles.xml part (wired into
# [Content_Types].xml + word/_rels/document.xml.rels automatically)
# containing an explicit paragraph style definition for "Normal".

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$normal = $d.Styles.Add("Normal", 1)
$normal.NameLocal = "Normal"

Write-Host ("Styles.Count=" + $d.Styles.Count)
Write-Host ("Normal style present: " + $normal.NameLocal)
